$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4315.4595
$ws.Range("I11").Value = 4315.4595
$ws.Range("K11").Value = 4315.4595
$ws.Range("M11").Value = -4175.4595
$ws.Range("H17").Value = 3327.1333
$ws.Range("J17").Value = 3327.1333
$ws.Range("L17").Value = 9981.3999
$ws.Range("N17").Value = -10317.3999
$ws.Range("H33").Value = 372.5
$ws.Range("I33").Value = 485
$ws.Range("J33").Value = 297.5
$ws.Range("K33").Value = 485
$ws.Range("L33").Value = 297.5
$ws.Range("M33").Value = -256
$ws.Range("N33").Value = -755.5
$ws.Range("H43").Value = 16667.166
$ws.Range("I43").Value = 22000.666
$ws.Range("J43").Value = 14889.333
$ws.Range("K43").Value = 22000.666
$ws.Range("L43").Value = 14889.333
$ws.Range("M43").Value = -21931.666
$ws.Range("N43").Value = -15027.333
$ws.Range("H47").Value = 57119.6
$ws.Range("I47").Value = 36533
$ws.Range("K47").Value = 36533
$ws.Range("M47").Value = -35561
$ws.Range("H51").Value = 10179.88
$ws.Range("I51").Value = 6388.8887
$ws.Range("K51").Value = 6388.8887
$ws.Range("M51").Value = -5904.8887
$ws.Range("H53").Value = 527.7143
$ws.Range("I53").Value = 311.1613
$ws.Range("K53").Value = 311.1613
$ws.Range("M53").Value = 325.8387
$ws.Range("H62").Value = 29636740
$ws.Range("I62").Value = 29636740
$ws.Range("K62").Value = 29636740
$ws.Range("M62").Value = -29636116
$ws.Range("H65").Value = 29636740
$ws.Range("I65").Value = 29636740
$ws.Range("K65").Value = 148183700
$ws.Range("M65").Value = -148180580
$ws.Range("H74").Value = 6928.9287
$ws.Range("I74").Value = 4800.6
$ws.Range("J74").Value = 8111.3335
$ws.Range("K74").Value = 4800.6
$ws.Range("L74").Value = 8111.3335
$ws.Range("M74").Value = -3864.6
$ws.Range("N74").Value = -9983.333500000001
$ws.Range("H77").Value = 6928.9287
$ws.Range("I77").Value = 4800.6
$ws.Range("J77").Value = 8111.3335
$ws.Range("K77").Value = 24003
$ws.Range("L77").Value = 40556.6675
$ws.Range("M77").Value = -19323
$ws.Range("N77").Value = -49916.6675
$ws.Range("H112").Value = 1762.3846
$ws.Range("J112").Value = 1809.25
$ws.Range("L112").Value = 5427.75
$ws.Range("N112").Value = -7643.75
$ws.Range("H132").Value = 998.1053000000001
$ws.Range("I132").Value = 998.1053000000001
$ws.Range("K132").Value = 2994.3159
$ws.Range("M132").Value = -464.3159000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 157090.5
$ws.Range("J130").Value = 157090.5
$ws.Range("L130").Value = 157090.5
$ws.Range("N130").Value = -167130.5
$ws.Range("H132").Value = 3333.5247
$ws.Range("I132").Value = 2910.0173
$ws.Range("K132").Value = 8730.0519
$ws.Range("M132").Value = -6200.0519

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 53333.332
$ws.Range("I138").Value = 20000
$ws.Range("J138").Value = 120000
$ws.Range("K138").Value = 20000
$ws.Range("L138").Value = 120000
$ws.Range("M138").Value = -14860
$ws.Range("N138").Value = -130280

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 53417.168
$ws.Range("J59").Value = 68466.336
$ws.Range("L59").Value = 68466.336
$ws.Range("N59").Value = -70756.336
$ws.Range("H60").Value = 19616.334
$ws.Range("J60").Value = 29299.666
$ws.Range("L60").Value = 29299.666
$ws.Range("N60").Value = -30321.666
$ws.Range("H97").Value = 59997.5
$ws.Range("J97").Value = 59997.5
$ws.Range("L97").Value = 59997.5
$ws.Range("N97").Value = -61979.5
$ws.Range("H100").Value = 45870
$ws.Range("J100").Value = 45870
$ws.Range("L100").Value = 45870
$ws.Range("N100").Value = -48034
$ws.Range("H134").Value = 2736
$ws.Range("I134").Value = 2335.5
$ws.Range("K134").Value = 7006.5
$ws.Range("M134").Value = -4471.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 74300.07000000001
$ws.Range("J37").Value = 74300.07000000001
$ws.Range("L37").Value = 222900.21
$ws.Range("N37").Value = -223124.21
$ws.Range("I47").Value = 166867.17
$ws.Range("J47").Value = 700
$ws.Range("K47").Value = 500601.51
$ws.Range("L47").Value = 2100
$ws.Range("M47").Value = -500170.51
$ws.Range("N47").Value = -2962
$ws.Range("H129").Value = 33334178
$ws.Range("J129").Value = 125001770
$ws.Range("L129").Value = 375005310
$ws.Range("N129").Value = -375015310
$ws.Range("H131").Value = 7756164.5
$ws.Range("J131").Value = 8554.421
$ws.Range("L131").Value = 25663.263
$ws.Range("N131").Value = -35743.263

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 204782.8
$ws.Range("I7").Value = 1000000
$ws.Range("J7").Value = 5978.5
$ws.Range("K7").Value = 1000000
$ws.Range("L7").Value = 5978.5
$ws.Range("M7").Value = -999888
$ws.Range("N7").Value = -6202.5
$ws.Range("H35").Value = 2291
$ws.Range("I35").Value = 2291
$ws.Range("K35").Value = 2291
$ws.Range("M35").Value = -1955
$ws.Range("H39").Value = 31174.25
$ws.Range("I39").Value = 28232.666
$ws.Range("K39").Value = 28232.666
$ws.Range("M39").Value = -27772.666
$ws.Range("H68").Value = 2585
$ws.Range("I68").Value = 2220
$ws.Range("J68").Value = 3497.5
$ws.Range("K68").Value = 2220
$ws.Range("L68").Value = 3497.5
$ws.Range("M68").Value = -1471
$ws.Range("N68").Value = -4995.5
$ws.Range("H71").Value = 2585
$ws.Range("I71").Value = 2220
$ws.Range("J71").Value = 3497.5
$ws.Range("K71").Value = 11100
$ws.Range("L71").Value = 17487.5
$ws.Range("M71").Value = -7356
$ws.Range("N71").Value = -24975.5
$ws.Range("H126").Value = 204782.8
$ws.Range("I126").Value = 1000000
$ws.Range("J126").Value = 5978.5
$ws.Range("K126").Value = 3000000
$ws.Range("L126").Value = 17935.5
$ws.Range("M126").Value = -2997530
$ws.Range("N126").Value = -22875.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 27750
$ws.Range("H54").Value = 153333
$ws.Range("I54").Value = 103333.336
$ws.Range("K54").Value = 103333.336
$ws.Range("M54").Value = -102813.336
$ws.Range("H81").Value = 11457.346
$ws.Range("I81").Value = 2586.5293
$ws.Range("J81").Value = 15425.868
$ws.Range("K81").Value = 5173.0586
$ws.Range("L81").Value = 30851.736
$ws.Range("M81").Value = -4112.0586
$ws.Range("N81").Value = -32973.736
$ws.Range("H84").Value = 11457.346
$ws.Range("I84").Value = 2586.5293
$ws.Range("J84").Value = 15425.868
$ws.Range("K84").Value = 25865.293
$ws.Range("L84").Value = 154258.68
$ws.Range("M84").Value = -20561.293
$ws.Range("N84").Value = -164866.68
$ws.Range("H126").Value = 3910.7334
$ws.Range("J126").Value = 4066.6667
$ws.Range("L126").Value = 12200.0001
$ws.Range("N126").Value = -17140.0001
$ws.Range("H136").Value = 2861.15
$ws.Range("I136").Value = 1367.7727
$ws.Range("J136").Value = 4686.3887
$ws.Range("K136").Value = 4103.3181
$ws.Range("L136").Value = 14059.1661
$ws.Range("M136").Value = -1553.3181
$ws.Range("N136").Value = -19159.1661
